$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (rows 5-102, the used range),
# shifting the existing quarterly data from D:K to F:M
$ws.Range("D5:E102").Insert(-4161)

# Copy number formatting/styles from column F (the old column D, now shifted)
# into the new D:E columns so the new quarters look like the rest of the table
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new D/E columns with the two most recent quarters of data
$ws.Range("D7").Value = 43462
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 452000
$ws.Range("E8").Value = 477600
$ws.Range("D9").Value = 341800
$ws.Range("E9").Value = 365800
$ws.Range("D10").Value = 110200
$ws.Range("E10").Value = 111800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 1400
$ws.Range("E14").Value = 600
$ws.Range("D15").Value = 8200
$ws.Range("E15").Value = 8000
$ws.Range("D17").Value = 406400
$ws.Range("E17").Value = 430900
$ws.Range("D18").Value = 45600
$ws.Range("E18").Value = 46700
$ws.Range("D20").Value = -10500
$ws.Range("E20").Value = -12600
$ws.Range("D21").Value = 53100
$ws.Range("E21").Value = 51800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 35100
$ws.Range("E23").Value = 34100
$ws.Range("D24").Value = 8200
$ws.Range("E24").Value = 1400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 26900
$ws.Range("E26").Value = 32700
$ws.Range("D27").Value = 26300
$ws.Range("E27").Value = 32700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 10500
$ws.Range("E32").Value = 12600
$ws.Range("D33").Value = 26300
$ws.Range("E33").Value = 32700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 26300
$ws.Range("E35").Value = 32700
$ws.Range("D38").Value = 43462
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 75900
$ws.Range("E41").Value = 126700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 294200
$ws.Range("E43").Value = 265100
$ws.Range("D44").Value = 226600
$ws.Range("E44").Value = 221800
$ws.Range("D45").Value = 29000
$ws.Range("E45").Value = 33600
$ws.Range("D46").Value = 625600
$ws.Range("E46").Value = 647100
$ws.Range("D47").Value = 2200
$ws.Range("E47").Value = 1600
$ws.Range("D48").Value = 241700
$ws.Range("E48").Value = 213100
$ws.Range("D49").Value = 474000
$ws.Range("E49").Value = 462000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1200
$ws.Range("E52").Value = 200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1344800
$ws.Range("E54").Value = 1324100
$ws.Range("D57").Value = 134100
$ws.Range("E57").Value = 156500
$ws.Range("D58").Value = 26600
$ws.Range("E58").Value = 26600
$ws.Range("D59").Value = 126900
$ws.Range("E59").Value = 89700
$ws.Range("D60").Value = 287500
$ws.Range("E60").Value = 272700
$ws.Range("D61").Value = 878100
$ws.Range("E61").Value = 877700
$ws.Range("D62").Value = 55000
$ws.Range("E62").Value = 51600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1220600
$ws.Range("E66").Value = 1202000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -314800
$ws.Range("E72").Value = -317400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 124200
$ws.Range("E76").Value = 122100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43462
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 26300
$ws.Range("E81").Value = 32700
$ws.Range("D83").Value = 18000
$ws.Range("E83").Value = 17600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 40300
$ws.Range("E89").Value = 25100
$ws.Range("D91").Value = -6900
$ws.Range("E91").Value = -12200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -64900
$ws.Range("E94").Value = -11900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -25200
$ws.Range("E100").Value = 5100
$ws.Range("D101").Value = -900
$ws.Range("E101").Value = -1200
$ws.Range("D102").Value = -50700
$ws.Range("E102").Value = 17100

# Apply data corrections identified in the restated figures for shifted columns
$ws.Range("I9").Value = 306300
$ws.Range("I10").Value = 89500
$ws.Range("I17").Value = 356800
$ws.Range("I18").Value = 39000
$ws.Range("I20").Value = -5900
$ws.Range("H24").Value = 7300
$ws.Range("H26").Value = 22400
$ws.Range("H27").Value = 22200
$ws.Range("H29").Value = 4800
$ws.Range("I32").Value = 5900
$ws.Range("H33").Value = 27000
$ws.Range("H35").Value = 27000
$ws.Range("H81").Value = 27000
